$wb = $excel.ActiveWorkbook

# --- "Data" sheet: append two new weekly WALCL observations ---
$data = $wb.Worksheets.Item("Data")

# Row 110 - copy formatting from the last existing data row (109) so the
# new date cell keeps the same date style/border/alignment as the rest of
# column A, then fill in the values.
$data.Cells.Item(109, 1).Copy()
$data.Cells.Item(110, 1).PasteSpecial(-4122)
$data.Cells.Item(110, 1).Value = 45231
$data.Cells.Item(110, 2).Value = 7866.664

# Row 111
$data.Cells.Item(109, 1).Copy()
$data.Cells.Item(111, 1).PasteSpecial(-4122)
$data.Cells.Item(111, 1).Value = 45238
$data.Cells.Item(111, 2).Value = 7860.691

$excel.CutCopyMode = $false

# --- "SeriesInfo" sheet: refresh FRED pull metadata ---
# These values must stay literal text (e.g. "2023-11-15"), not get
# auto-converted to date serials by a plain .Value assignment, and must not
# pick up a new number-format style. Writing them as a quoted-string formula
# and then collapsing the formula to its static value via PasteSpecial
# (values only) keeps the cell a plain text cell with the default style.
$info = $wb.Worksheets.Item("SeriesInfo")

function Set-TextValue($range, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

Set-TextValue $info.Range("B3") "2023-11-15"
Set-TextValue $info.Range("B4") "2023-11-15"
Set-TextValue $info.Range("B7") "2023-11-08"
Set-TextValue $info.Range("B14") "2023-11-09 15:37:01-06"

$excel.CutCopyMode = $false
